$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tabla1")

# --- Add a brand-new product row at the bottom of the table (row 228) ---
# This must happen FIRST so its new shared-string value ("intercyprus...")
# is appended to the shared-strings table before the other new values below.
$newRow = $tbl.ListRows.Add()

$ws.Range("H227").Copy()
$ws.Range("H228").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E228").Value = "Pegatina"
$ws.Range("F228").Formula = "=+Tabla1[[#This Row],[CODIGO]]"
$ws.Range("I228").Value = "https://www.intercyprus.com/cdn/shop/files/aaf729e2-cc35-4d14-96d6-b381f5aa6252.jpg"

# --- Fill in missing IMAGEN (and one MODAL) links for existing rows, in row order ---
$ws.Range("I124").Value = "https://img.joomcdn.net/80c1a60906182457d02c893b1c74dbbca0ebe7ed_original.jpeg"
$ws.Range("I126").Value = "https://img.joomcdn.net/fe719753317e67740a5a055ad1d889b2db6c5391_original.jpeg"

$ws.Range("I132").Value = "https://coolboxpe.vtexassets.com/arquivos/ids/277465-800-800?v=638219163682200000.jpg"
$ws.Range("J132").Value = "https://plazavea.vteximg.com.br/arquivos/ids/16755450-418-418/image-4926ee5b681b412b9fd891f7bee418af.jpg"

$ws.Range("I133").Value = "https://dojiw2m9tvv09.cloudfront.net/76107/product/proco-14x-e13455.jpg"
$ws.Range("I134").Value = "https://http2.mlstatic.com/D_NQ_NP_841592-MLU54718848989_032023-O.jpg"
$ws.Range("I136").Value = "https://m.media-amazon.com/images/I/71k+JJfLdDL._UY350_.jpg"

$ws.Range("I158").Value = "https://s.alicdn.com/@sc04/kf/H5e5f3b18d84940178b3e2ed9c7b0f862s.jpg"
$ws.Range("I159").Value = "https://tenda24.pe/cdn/shop/files/S1d2c150f1f864d0badec822bbfcfd0c30.webp"
$ws.Range("I160").Value = "https://http2.mlstatic.com/D_Q_NP_735491-MLU73805295276_012024-O.jpg"
$ws.Range("I162").Value = "https://img.joomcdn.net/4011e499955db4dacb10772d332db88ee38c5a17_original.jpeg"
$ws.Range("I163").Value = "https://ae01.alicdn.com/kf/S5ed3acb0573b4631b7cb40a2eecbe120g.jpg"
$ws.Range("I165").Value = "https://ae01.alicdn.com/kf/S1c9be3d8efa84444ab3361d982856993e.jpg"

$ws.Range("I168").Value = "https://plazavea.vteximg.com.br/arquivos/ids/25993859-418-418/image-6f4072c3b49c4143bf7462bb524e3d01.jpg"
$ws.Range("I169").Value = "https://m.media-amazon.com/images/I/71uQA-j8elL._UF894,1000_QL80_.jpg"

$ws.Range("I175").Value = "https://m.media-amazon.com/images/I/81puHsA9KAL._UY350_.jpg"

$ws.Range("I182").Value = "https://m.media-amazon.com/images/I/71-9vBvpgDL._UF894,1000_QL80_.jpg"
$ws.Range("I183").Value = "https://http2.mlstatic.com/D_NQ_NP_791879-MLU71266333911_082023-O.jpg"

$ws.Range("H192").Value = 28
$ws.Range("I192").Value = "https://adhek-peru.com/cdn/shop/products/FIGURA_BROLY_DRAGON_BALLZ_26114_BANDAI_a_2048x.jpg"

$ws.Range("I203").Value = "https://i5.walmartimages.com/asr/4408004d-64f3-43a7-b4ab-6e0be6259cd4.7bef667ddcb77ffe921d8ad96d967b84.jpeg"
$ws.Range("I205").Value = "https://ae01.alicdn.com/kf/S2fda666123c34aa9ac4437dfd66ec1735.jpg"
$ws.Range("I206").Value = "https://ae01.alicdn.com/kf/Sf158e6ce530a4cfb9e74bdd4fe91eeeej.jpg"

# --- Update the view: active selection (also scrolls the window so the
#     selected cell becomes the new top-left visible cell, matching the
#     saved sheetView's topLeftCell/selection state) ---
$ws.Activate()
$ws.Range("I201").Select()
